# Junction_Flooding_173.xlsx edit:
#  - row 5 (B5:AH5) values re-exported with "custom accuracy" (2 decimal places
#    instead of 3) per commit message "custom accuracy + 데이터 1000개"
#  - row 6 removed (dataset trimmed down)
#  - column L width narrowed from 8 to 7 character units

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 5 values (B5:AH5) to the new rounded readings ---
$values = @(16.95, 12.83, 0.74, 36.41, 30.5, 13.09, 50.47, 20.11, 9.13, 13.83, `
            14.75, 15.46, 4.2, 13.04, 18.77, 10.87, 0.45, 0.37, 192.79, 36.54, `
            11.99, 24.72, 13.26, 1.73, 25.24, 10.41, 9.92, 10.9, 15.52, 0, `
            45.41, 7.08, 15.02)

$arr = New-Object 'object[,]' 1,33
for ($i = 0; $i -lt 33; $i++) {
    $arr[0,$i] = $values[$i]
}
$ws.Range("B5:AH5").Value = $arr

# --- Remove row 6 entirely (shifts dimension from A1:AH6 to A1:AH5) ---
$ws.Rows.Item(6).Delete()

# --- Narrow column L (12) from width 8 to width 7 ---
# ColumnWidth uses character units that serialize as (raw_width - 5/6); solve
# backwards so the stored <col width="..."> ends up exactly at 7.
$ws.Columns.Item(12).ColumnWidth = 7 - (5/6)
